# Replace the old per-build rsid / bookmark-id stamps with the new ones
# (tool moved from version 3.1.0 to 3.1.1, which re-generated these
# internal identifiers). Content/text itself is unchanged.
#
# Because bookmark/rsid attributes are not exposed as settable
# properties on the Word object model, each affected paragraph is
# rebuilt in place (same runs/content/formatting) via Range.InsertXML
# with the literal OOXML carrying the new attribute values.

$d = $word.ActiveDocument

$oldRsid = "59E62E63122EB0AF6716498485023C18"
$newRsid = "72FE6F9E06FBBB0C64CBF15D93BB830F"
$oldBookmarkId = "116620715360368245400264450396621043482"
$newBookmarkId = "20613205666473705645396298935817703149"

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph "Test link before bookmark : <REF field>" ---
$para1 = '<w:body><w:p w:rsidP="009168BC" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Test link before bookmark : </w:t></w:r><w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="' + $newRsid + '"><w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText></w:r><w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="' + $newRsid + '"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body>'

$r1 = $d.Paragraphs(2).Range
$r1.InsertXML($pkgOpen + $para1 + $pkgClose)

# --- Paragraph "Test bookmark : bookmarked content" ---
$para2 = '<w:body><w:p w:rsidP="009168BC" w:rsidR="00C52979" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t>Test</w:t></w:r><w:r w:rsidR="00C52979"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>bookmark</w:t></w:r><w:r w:rsidR="00C52979"><w:t xml:space="preserve"> : </w:t></w:r><w:bookmarkStart w:name="bookmark1" w:id="' + $newBookmarkId + '"/><w:r><w:t>bookmarked content</w:t></w:r><w:bookmarkEnd w:id="' + $newBookmarkId + '"/></w:p></w:body>'

$r2 = $d.Paragraphs(3).Range
$r2.InsertXML($pkgOpen + $para2 + $pkgClose)

# --- Paragraph "Test link after bookmark : <REF field>" ---
$para3 = '<w:body><w:p w:rsidP="00E02A2B" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Test link after bookmark : </w:t></w:r><w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="' + $newRsid + '"><w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText></w:r><w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="' + $newRsid + '"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidR="00D0546C"><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body>'

$r3 = $d.Paragraphs(4).Range
$r3.InsertXML($pkgOpen + $para3 + $pkgClose)
